$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Cells.Item(19, 2).Value = 44138
$ws.Cells.Item(19, 3).Value = 0.58333333333333337
$ws.Cells.Item(19, 4).Value = 0.95833333333333337

$ws.Cells.Item(20, 2).Value = 44139
$ws.Cells.Item(20, 3).Value = 0.5625
$ws.Cells.Item(20, 4).Value = 0.875

$ws.Cells.Item(21, 2).Value = 44139
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 4).Value = 0.0625

$ws.Cells.Item(22, 2).Value = 44140
$ws.Cells.Item(22, 3).Value = 0.64583333333333337
$ws.Cells.Item(22, 4).Value = 0.6875

$ws.Cells.Item(23, 2).Value = 44141
$ws.Cells.Item(23, 3).Value = 0.47916666666666669
$ws.Cells.Item(23, 4).Value = 0.57638888888888895

$ws.Cells.Item(24, 2).Value = 44145
$ws.Cells.Item(24, 3).Value = 0.41666666666666669
$ws.Cells.Item(24, 4).Value = 0.4513888888888889

$ws.Cells.Item(25, 2).Value = 44145
$ws.Cells.Item(25, 3).Value = 0.875
$ws.Cells.Item(25, 4).Value = 0.9375

$ws.Cells.Item(26, 2).Value = 44146
$ws.Cells.Item(26, 3).Value = 0.041666666666666664
$ws.Cells.Item(26, 4).Value = 0.22916666666666666

$excel.Calculate()

$ws.Range("D26").Select()
